$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.240655899047852
$ws.Range("B1").Value = 2.363288879394531
$ws.Range("C1").Value = 3.809030771255493
$ws.Range("D1").Value = 3.114845037460327
$ws.Range("E1").Value = 1.284807920455933
